# Results_IA_01-02.xlsx -- "various updates onP Part 1.3 and Part 1.4"
#
# Adds a second (Negative-word) table in C:D mirroring the existing
# Positive-word table in A:B, a small Train/Validation/Test accuracy
# summary in G2:I3, and a short "Correct classified / Falsely classified"
# example-sentences block under the existing data (rows 15-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Negative-word table header (C1 merged across C1:D1, centered like A1:B1)
# ---------------------------------------------------------------------
$ws.Range("C1:D1").Merge()
$ws.Range("C1").Value = "Negative"
$ws.Range("C1:D1").HorizontalAlignment = -4108   # xlCenter

# Column header row (C2/D2 already hold "Word"/"Weights" per original diff context)
$ws.Range("D2").Value = "Weights"

# ---------------------------------------------------------------------
# 2) Negative-word weight values, D3:D12 (same 0.0000 number format as B3:B12)
# ---------------------------------------------------------------------
$ws.Range("D3").NumberFormat = "0.0000"
$ws.Range("D3").Value = -7.3593803962037798
$ws.Range("D4").NumberFormat = "0.0000"
$ws.Range("D4").Value = -6.9300100000000002
$ws.Range("D5").NumberFormat = "0.0000"
$ws.Range("D5").Value = -6.1804800000000002
$ws.Range("D6").NumberFormat = "0.0000"
$ws.Range("D6").Value = -5.7775400000000001
$ws.Range("D7").NumberFormat = "0.0000"
$ws.Range("D7").Value = -5.08507
$ws.Range("D8").NumberFormat = "0.0000"
$ws.Range("D8").Value = -5.0396939999999999
$ws.Range("D9").NumberFormat = "0.0000"
$ws.Range("D9").Value = -4.8361999999999998
$ws.Range("D10").NumberFormat = "0.0000"
$ws.Range("D10").Value = -4.5254000000000003
$ws.Range("D11").NumberFormat = "0.0000"
$ws.Range("D11").Value = -4.3102
$ws.Range("D12").NumberFormat = "0.0000"
$ws.Range("D12").Value = -5.2602000000000002
$ws.Range("D13").NumberFormat = "0.0000"

$ws.Columns("D").ColumnWidth = 11.6640625

# ---------------------------------------------------------------------
# 3) Train/Validation/Test accuracy summary, G2:I3
# ---------------------------------------------------------------------
$ws.Range("G2").Value = "Training"
$ws.Range("H2").Value = "Validation"
$ws.Range("I2").Value = "Test"

$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0.96499999999999997
$ws.Range("I3").Value = 0.96

# ---------------------------------------------------------------------
# 4) Correctly / falsely classified example sentences, rows 15-21
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "Correct classified:"

$ws.Range("A16").Font.Size = 10
$ws.Range("A16").Font.Name = "Helvetica"
$ws.Range("A16").Value = " before boarding , the terrible customer service people at the gate made people check their bags even though there was an excessive amount of space in overhead compartments . with their terrible rate of baggage loss , this has ruined my flight and others ."

$ws.Range("A17").Font.Size = 10
$ws.Range("A17").Font.Name = "Helvetica"
$ws.Range("A17").Value = "london to delhi . an excellent service and experience . this was my first time travelling with ai and i was amazed with service from ground staff to onboard . i'm looking forward flying with ai again in future ."

$ws.Range("A19").Value = "Flasly Classified"

$ws.Range("A20").Value = " london heathrow to algiers on an old 767 . cabin crew are simply not professional food wasn't appealing and toilets were filthy ."

$ws.Range("A21").Value = "terrible service . mean and unkind employees . left me waiting 46 minutes to find a wheelchair . doesn’t have enough staff to operate wheelchair or check everyone in . horrible food ."

# ---------------------------------------------------------------------
# 5) Selection ends on A21, matching the saved cursor position
# ---------------------------------------------------------------------
$ws.Range("A21").Select()
